$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NroPoliza and FechaSiniestro values on row 2
$ws.Range("E2").Value = "'11111003260"
$ws.Range("G2").Value = "'20/05/2021"

# Update the active selection on the sheet
$ws.Range("E3").Select()
